$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Croatia HNL")

$rows = @(
    @{ Row=121; A=119; B=6787894; C="Croatia HNL"; D="Croatia HNL"; E=45353.45833333334; F="Slaven Belupo"; G="NK Lokomotiva Zagreb"; K=2.9;  L=3.2;  M=2.3;   N=2.875; O=3;    P=2.45; Q=0;    R=2.05;  S=1.8;   T=2.25; U=1.9;   V=1.95;  W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=122; A=120; B=6788928; C="Croatia HNL"; D="Croatia HNL"; E=45353.54861111111; F="Hajduk Split"; G="Istra 1961";            K=1.3;  L=5;    M=9.5;   N=1.3;   O=5;    P=9.5;  Q=-1.5; R=1.95;  S=1.9;   T=2.5;  U=1.9;   V=1.95;  W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=123; A=121; B=6788926; C="Croatia HNL"; D="Croatia HNL"; E=45354.45833333334; F="HNK Rijeka";   G="NK Varazdin";            K=1.25; L=4.75; M=11;    N=1.25;  O=4.75; P=11;   Q=-1.5; R=1.9;   S=1.95;  T=2.5;  U=1.825; V=2.025; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=124; A=122; B=6788927; C="Croatia HNL"; D="Croatia HNL"; E=45354.5625;        F="NK Osijek";    G="Dinamo Zagreb";          K=3.6;  L=3.3;  M=1.909; N=4.333; O=3.3;  P=1.8;  Q=0.5;  R=2.05;  S=1.8;   T=2.25; U=1.825; V=2.025; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=125; A=123; B=6769303; C="Croatia HNL"; D="Croatia HNL"; E=45355.54166666666; F="NK Rudes";     G="HNK Gorica";             K=3.3;  L=3.2;  M=2.1;   N=4;     O=3.1;  P=1.909; Q=0.5; R=1.825; S=2.025; T=2.25; U=1.925; V=1.925; W=0; X=0; Y=0; Z=0; AA=0 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Copy cell formatting (styles) from the last pre-existing data row (120)
    # so the new rows inherit the same "id" and "Date" column formatting.
    $ws.Range("A120").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("E120").Copy() | Out-Null
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D

    $ws.Cells.Item($row, 5).Value = $r.E

    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G

    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
    $ws.Cells.Item($row, 23).Value = $r.W
    $ws.Cells.Item($row, 24).Value = $r.X
    $ws.Cells.Item($row, 25).Value = $r.Y
    $ws.Cells.Item($row, 26).Value = $r.Z
    $ws.Cells.Item($row, 27).Value = $r.AA
}
